$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row: insert department_id column before is_active ---
$ws.Range("F1").Value = "department_id"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "is_active"

# --- Fix existing data rows (2-7): plan_no *111, plan_name fix, department_id, is_active moved to G ---
$ws.Range("B2").Value = 111
$ws.Range("C2").Value = "بكالوريوس هندسة نظم الحاسوب"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("B3").Value = 222
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("B4").Value = 333
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1

$ws.Range("B5").Value = 444
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 1

$ws.Range("B6").Value = 555
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 1

$ws.Range("B7").Value = 666
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 1

# --- New rows ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 777
$ws.Range("C9").Value = "دبلوم التحكم الصناعي"
$ws.Range("D9").Value = 2020
$ws.Range("E9").Value = 73
$ws.Range("F9").Value = "CSE"
$ws.Range("G9").Value = 1

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 888
$ws.Range("C10").Value = "دبلوم محاسبة"
$ws.Range("D10").Value = 2020
$ws.Range("E10").Value = 73
$ws.Range("F10").Value = "العلوم الإدارية والمالية"
$ws.Range("G10").Value = 1

$ws.Range("A11").Value = "يعمري انتا ما احلاك"
$ws.Range("A11:F11").HorizontalAlignment = -4108
$ws.Range("A11:F11").Merge()

$ws.Range("A12").Value = 9
$ws.Range("B12").Value = 999
$ws.Range("C12").Value = "إعلام رقمي"
$ws.Range("D12").Value = 2020
$ws.Range("E12").Value = 73
$ws.Range("F12").Value = "العلوم التطبيقية"
$ws.Range("G12").Value = 1

$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 999
$ws.Range("C13").Value = "اعلام رقمي"
$ws.Range("D13").Value = 2020
$ws.Range("E13").Value = 73
$ws.Range("F13").Value = "العلوم التطبيقية"
$ws.Range("G13").Value = 1

$ws.Range("A14").Value = 11
$ws.Range("B14").Value = 666
$ws.Range("C14").Value = "دبلوم التحكم الصناعي"
$ws.Range("D14").Value = 2020
$ws.Range("E14").Value = 73
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1

# --- Column widths (COM ColumnWidth is quantized to 1/7-character pixel
#     steps; these are the closest achievable values to the target OOXML
#     widths of 16 and 13.25) ---
$ws.Range("F1").ColumnWidth = 15.2857142857143
$ws.Range("G1").ColumnWidth = 12.5714285714286

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ---
$ws.Range("A13").Select()
